$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 15.28448560880142)
    3 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    5 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    6 = @(0.04271373187048222, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.115428400803308)
    7 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523)
    8 = @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.406728370586922)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
